# 状態図.xlsx — "up to what before 考察"
#
# 1) Worksheet: add a 4th ("ratio of Sn" in D1 of the 温度（固相線）column
#    header), re-label the headers in English for columns A-C, moving the
#    original Japanese header out to column D.
# 2) Chart: give the chart a title ("状態図") and name the two series
#    ("固相線" / "液相線") so the legend reflects them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet headers -----------------------------------------------
# Before:  A1="Snの割合"            B1="温度（固相線）"
# After:   A1="ratio of Sn"         B1="solid phase line"
#          C1="liquid phase line"   D1="温度（固相線）"
# (Move the original Japanese temperature header to D1 first, then
#  overwrite A1:C1 with the new English headers.)
$ws.Range("D1").Value = "温度（固相線）"
$ws.Range("B1").Value = "solid phase line"
$ws.Range("C1").Value = "liquid phase line"
$ws.Range("A1").Value = "ratio of Sn"

# Move the active selection off the stale "K21" left over from the
# author's last session back onto the top-left cell.
[void]$ws.Range("A1").Select()

# --- Chart title + series names ---------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$chart.HasTitle = $true
$chart.ChartTitle.Text = "状態図"

$chart.SeriesCollection(1).Name = "固相線"
$chart.SeriesCollection(2).Name = "液相線"
